# Update "想去人数" (column F) counts on the 展览 / 演出 / 全部类型 sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 4704
$ws1.Range("F3").Value  = 1864
$ws1.Range("F6").Value  = 3157
$ws1.Range("F7").Value  = 584
$ws1.Range("F8").Value  = 594
$ws1.Range("F13").Value = 402
$ws1.Range("F16").Value = 1371
$ws1.Range("F18").Value = 1639
$ws1.Range("F21").Value = 615
$ws1.Range("F23").Value = 47
$ws1.Range("F24").Value = 540
$ws1.Range("F30").Value = 38
$ws1.Range("F32").Value = 3971
$ws1.Range("F34").Value = 779
$ws1.Range("F35").Value = 83
$ws1.Range("F36").Value = 1344
$ws1.Range("F38").Value = 1882

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 27
$ws2.Range("F3").Value = 55

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 4704
$ws4.Range("F3").Value  = 1864
$ws4.Range("F6").Value  = 3157
$ws4.Range("F7").Value  = 584
$ws4.Range("F8").Value  = 594
$ws4.Range("F13").Value = 27
$ws4.Range("F14").Value = 402
$ws4.Range("F17").Value = 1371
$ws4.Range("F19").Value = 1639
$ws4.Range("F22").Value = 615
$ws4.Range("F24").Value = 47
$ws4.Range("F25").Value = 540
$ws4.Range("F31").Value = 38
$ws4.Range("F33").Value = 3971
$ws4.Range("F34").Value = 55
$ws4.Range("F37").Value = 779
$ws4.Range("F38").Value = 83
$ws4.Range("F39").Value = 1344
$ws4.Range("F41").Value = 1882
